$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 27
$ws.Range("A27").Value = 1747906147
$ws.Range("B27").Value = 'update'
$ws.Range("C27").Value = 'variable'
$ws.Range("D27").Value = 'ser_pub_loc___variable_12'
$ws.Range("F27").Value = 'source_var_ids'
$ws.Range("H27").Value = 'ser_pub_loc___variable_10, ser_pub_loc___variable_11'

# Row 28
$ws.Range("A28").Value = 1747906281
$ws.Range("B28").Value = 'add'
$ws.Range("C28").Value = 'config'
$ws.Range("D28").Value = 'alias_3'

# Row 29
$ws.Range("A29").Value = 1747906551
$ws.Range("B29").Value = 'update'
$ws.Range("C29").Value = 'variable'
$ws.Range("D29").Value = 'ser_pub_loc___variable_12'
$ws.Range("F29").Value = 'source_var_ids'
$ws.Range("G29").Value = 'ser_pub_loc___variable_10, ser_pub_loc___variable_11'

# Row 30
$ws.Range("A30").Value = 1747906551
$ws.Range("B30").Value = 'update'
$ws.Range("C30").Value = 'variable'
$ws.Range("D30").Value = 'ser_pub_loc___variable_12'
$ws.Range("F30").Value = 'sourceVar_ids'
$ws.Range("H30").Value = 'ser_pub_loc___variable_10, ser_pub_loc___variable_11'

# Row 31
$ws.Range("A31").Value = 1747906566
$ws.Range("B31").Value = 'update'
$ws.Range("C31").Value = 'config'
$ws.Range("D31").Value = 'alias_3'
$ws.Range("F31").Value = 'value'
$ws.Range("G31").Value = 'variable : source_var'
$ws.Range("H31").Value = 'variable : sourceVar'

# Row 32
$ws.Range("A32").Value = 1747907576
$ws.Range("B32").Value = 'add'
$ws.Range("C32").Value = 'config'
$ws.Range("D32").Value = 'test'

# Row 33
$ws.Range("A33").Value = 1747907576
$ws.Range("B33").Value = 'delete'
$ws.Range("C33").Value = 'config'
$ws.Range("D33").Value = 'alias_3'

# Row 34
$ws.Range("A34").Value = 1747907625
$ws.Range("B34").Value = 'add'
$ws.Range("C34").Value = 'config'
$ws.Range("D34").Value = 'alias_3'

# Row 35
$ws.Range("A35").Value = 1747907625
$ws.Range("B35").Value = 'delete'
$ws.Range("C35").Value = 'config'
$ws.Range("D35").Value = 'test'

# Row 36
$ws.Range("A36").Value = 1747913177
$ws.Range("B36").Value = 'delete'
$ws.Range("C36").Value = 'config'
$ws.Range("D36").Value = 'alias_3'

# Row 37
$ws.Range("A37").Value = 1747913221
$ws.Range("B37").Value = 'update'
$ws.Range("C37").Value = 'variable'
$ws.Range("D37").Value = 'ser_pub_loc___variable_12'
$ws.Range("F37").Value = 'sourceVar_ids'
$ws.Range("G37").Value = 'ser_pub_loc___variable_10, ser_pub_loc___variable_11'

# Row 38
$ws.Range("A38").Value = 1747913221
$ws.Range("B38").Value = 'update'
$ws.Range("C38").Value = 'variable'
$ws.Range("D38").Value = 'ser_pub_loc___variable_12'
$ws.Range("F38").Value = 'source_ids'
$ws.Range("H38").Value = 'ser_pub_loc___variable_10, ser_pub_loc___variable_11'

# Row 39
$ws.Range("A39").Value = 1747913314
$ws.Range("B39").Value = 'update'
$ws.Range("C39").Value = 'variable'
$ws.Range("D39").Value = 'ser_pub_loc___variable_12'
$ws.Range("F39").Value = 'source_ids'
$ws.Range("G39").Value = 'ser_pub_loc___variable_10, ser_pub_loc___variable_11'

# Row 40
$ws.Range("A40").Value = 1747913314
$ws.Range("B40").Value = 'update'
$ws.Range("C40").Value = 'variable'
$ws.Range("D40").Value = 'ser_pub_loc___variable_12'
$ws.Range("F40").Value = 'sourceVar_ids'
$ws.Range("H40").Value = 'ser_pub_loc___variable_10, ser_pub_loc___variable_11'

# Row 41
$ws.Range("A41").Value = 1747918633
$ws.Range("B41").Value = 'update'
$ws.Range("C41").Value = 'variable'
$ws.Range("D41").Value = 'ser_pub_loc___variable_13'
$ws.Range("F41").Value = 'sourceVar_ids'
$ws.Range("H41").Value = 'ser_pub_loc___variable_12'

# Row 42
$ws.Range("A42").Value = 1747921769
$ws.Range("B42").Value = 'update'
$ws.Range("C42").Value = 'variable'
$ws.Range("D42").Value = 'ser_pub_loc___variable_14'
$ws.Range("F42").Value = 'sourceVar_ids'
$ws.Range("H42").Value = 'accident_route___variable_7, dep_sante___variable_7'

# Row 43
$ws.Range("A43").Value = 1747921769
$ws.Range("B43").Value = 'update'
$ws.Range("C43").Value = 'variable'
$ws.Range("D43").Value = 'ser_pub_loc___variable_15'
$ws.Range("F43").Value = 'sourceVar_ids'
$ws.Range("H43").Value = 'accident_route___variable_7, dep_sante___variable_7'

# Row 44
$ws.Range("A44").Value = 1747921769
$ws.Range("B44").Value = 'update'
$ws.Range("C44").Value = 'variable'
$ws.Range("D44").Value = 'dep_sante___variable_3'
$ws.Range("F44").Value = 'sourceVar_ids'
$ws.Range("H44").Value = 'accident_route___variable_7'

